$wb = $excel.ActiveWorkbook

# --- Sheet: Schedule ---
$schedule = $wb.Worksheets.Item("Schedule")

$schedule.Range("E2").Value = 365.2350682500001
$schedule.Range("F2").Value = 8.05191949404762
$schedule.Range("E3").Value = 430.938183
$schedule.Range("F3").Value = 28.50120257936508

# --- Sheet: Detailed ---
$detailed = $wb.Worksheets.Item("Detailed")

$detailed.Range("B10").Value = 73.2
$detailed.Range("B11").Value = 73.2
$detailed.Range("B12").Value = 77.94
$detailed.Range("C12").Value = "historical"
$detailed.Range("B13").Value = 80.97543
$detailed.Range("C13").Value = "historical"
$detailed.Range("B14").Value = 76.57747999999999
$detailed.Range("B15").Value = 58.85679
$detailed.Range("B18").Value = -5.50985
$detailed.Range("B19").Value = -6.41051
$detailed.Range("B20").Value = -7.86387
$detailed.Range("B21").Value = -7.75171
$detailed.Range("B22").Value = -6.93648
$detailed.Range("B23").Value = -7.48833
$detailed.Range("B24").Value = -7.58381
$detailed.Range("B25").Value = -6.57582
$detailed.Range("B26").Value = -6.37274
$detailed.Range("B27").Value = -6.6969
$detailed.Range("B28").Value = -7.86766
$detailed.Range("B29").Value = -6.10909
$detailed.Range("B30").Value = -5.50985
$detailed.Range("B31").Value = -2.54301
$detailed.Range("B34").Value = -9.78218
$detailed.Range("B35").Value = -12.01
$detailed.Range("B37").Value = -7.41943
$detailed.Range("B38").Value = -0.42828
$detailed.Range("B39").Value = 3.06245
$detailed.Range("B40").Value = 18.95808
$detailed.Range("B43").Value = 53.24602
$detailed.Range("B44").Value = 57.04367
$detailed.Range("B46").Value = 54.13155
$detailed.Range("B49").Value = 52.66655
